# Scheduled runner refresh: re-pull Universalis market prices and
# recompute Leve profit columns (H:N) for the affected leves on each
# Disciple of the Hand job sheet.
$wb = $excel.ActiveWorkbook

# ALC!row19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 556.5
$ws.Range("J19").Value = 558.4
$ws.Range("L19").Value = 558.4
$ws.Range("N19").Value = -908.4

# ALC!row51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7647.3335
$ws.Range("J51").Value = 7098.75
$ws.Range("L51").Value = 7098.75
$ws.Range("N51").Value = -8066.75

# ALC!row98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 33017.727
$ws.Range("I98").Value = 47246.332
$ws.Range("J98").Value = 15943.4
$ws.Range("K98").Value = 47246.332
$ws.Range("L98").Value = 15943.4
$ws.Range("M98").Value = -45748.332
$ws.Range("N98").Value = -18939.4

# ALC!row112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2801.5
$ws.Range("I112").Value = 1350
$ws.Range("J112").Value = 2962.7778
$ws.Range("K112").Value = 4050
$ws.Range("L112").Value = 8888.3334
$ws.Range("M112").Value = -2942
$ws.Range("N112").Value = -11104.3334

# ALC!row122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 33017.727
$ws.Range("I122").Value = 47246.332
$ws.Range("J122").Value = 15943.4
$ws.Range("K122").Value = 141738.996
$ws.Range("L122").Value = 47830.2
$ws.Range("M122").Value = -139288.996
$ws.Range("N122").Value = -52730.2

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 8116.7744
$ws.Range("I137").Value = 10863.904
$ws.Range("J137").Value = 2347.8
$ws.Range("K137").Value = 32591.712
$ws.Range("L137").Value = 7043.400000000001
$ws.Range("M137").Value = -30041.712
$ws.Range("N137").Value = -12143.4

# ARM!row63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1655.6666
$ws.Range("I63").Value = 1655.6666
$ws.Range("K63").Value = 1655.6666
$ws.Range("M63").Value = -969.6666

# ARM!row66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1655.6666
$ws.Range("I66").Value = 1655.6666
$ws.Range("K66").Value = 8278.333000000001
$ws.Range("M66").Value = -4846.333000000001

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3110.8
$ws.Range("J74").Value = 3841.889
$ws.Range("L74").Value = 3841.889
$ws.Range("N74").Value = -5589.889

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3110.8
$ws.Range("J77").Value = 3841.889
$ws.Range("L77").Value = 19209.445
$ws.Range("N77").Value = -27945.445

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2963.257
$ws.Range("I132").Value = 2395.8262
$ws.Range("J132").Value = 4050.8333
$ws.Range("K132").Value = 7187.4786
$ws.Range("L132").Value = 12152.4999
$ws.Range("M132").Value = -4657.4786
$ws.Range("N132").Value = -17212.4999

# BSM!row20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3085.6667
$ws.Range("I20").Value = 1903.2858
$ws.Range("J20").Value = 5450.4287
$ws.Range("K20").Value = 1903.2858
$ws.Range("L20").Value = 5450.4287
$ws.Range("M20").Value = -1656.2858
$ws.Range("N20").Value = -5944.4287

# BSM!row22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 78
$ws.Range("I22").Value = 78
$ws.Range("K22").Value = 78
$ws.Range("M22").Value = 95

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7092.857
$ws.Range("I134").Value = 7403.1
$ws.Range("K134").Value = 22209.3
$ws.Range("M134").Value = -19674.3

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8572.791999999999
$ws.Range("I31").Value = 10826.6
$ws.Range("J31").Value = 4816.4443
$ws.Range("K31").Value = 10826.6
$ws.Range("L31").Value = 4816.4443
$ws.Range("M31").Value = -10531.6
$ws.Range("N31").Value = -5406.4443

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 8572.791999999999
$ws.Range("I34").Value = 10826.6
$ws.Range("J34").Value = 4816.4443
$ws.Range("K34").Value = 10826.6
$ws.Range("L34").Value = 4816.4443
$ws.Range("M34").Value = -10624.6
$ws.Range("N34").Value = -5220.4443

# CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2628.5925
$ws.Range("I58").Value = 2706.0908
$ws.Range("K58").Value = 2706.0908
$ws.Range("M58").Value = -2503.0908

# CRP!row107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J107").Value = 915.25
$ws.Range("L107").Value = 915.25
$ws.Range("N107").Value = -4755.25

# CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2628.5925
$ws.Range("I136").Value = 2706.0908
$ws.Range("K136").Value = 8118.2724
$ws.Range("M136").Value = -5568.2724

# CUL!row7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 61.4
$ws.Range("I7").Value = 27.5
$ws.Range("K7").Value = 82.5
$ws.Range("M7").Value = 29.5

# CUL!row12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 135.13333
$ws.Range("J12").Value = 32.444443
$ws.Range("L12").Value = 97.33332899999999
$ws.Range("N12").Value = -443.333329

# GSM!row35
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").ClearContents()

# GSM!row70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6180289
$ws.Range("I70").Value = 8553033
$ws.Range("K70").Value = 8553033
$ws.Range("M70").Value = -8552763

# GSM!row73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6180289
$ws.Range("I73").Value = 8553033
$ws.Range("K73").Value = 8553033
$ws.Range("M73").Value = -8552097

# GSM!row123
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 38399.8
$ws.Range("J123").Value = 38399.8
$ws.Range("L123").Value = 38399.8
$ws.Range("N123").Value = -43299.8

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2779.48
$ws.Range("I132").Value = 2681.2273
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 8043.6819
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -5513.6819
$ws.Range("N132").Value = -15560

# LTW!row22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 28996.715
$ws.Range("I22").Value = 40358
$ws.Range("J22").Value = 593.5
$ws.Range("K22").Value = 40358
$ws.Range("L22").Value = 593.5
$ws.Range("M22").Value = -40063
$ws.Range("N22").Value = -1183.5

# LTW!row27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 28996.715
$ws.Range("I27").Value = 40358
$ws.Range("J27").Value = 593.5
$ws.Range("K27").Value = 40358
$ws.Range("L27").Value = 593.5
$ws.Range("M27").Value = -40251
$ws.Range("N27").Value = -807.5

# LTW!row98
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 355000
$ws.Range("J98").Value = 355000
$ws.Range("L98").Value = 355000
$ws.Range("N98").Value = -360990

# LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3646.2666
$ws.Range("J136").Value = 6420.5
$ws.Range("L136").Value = 19261.5
$ws.Range("N136").Value = -24361.5

# LTW!row140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 88997.39999999999
$ws.Range("J140").Value = 88997.39999999999
$ws.Range("L140").Value = 88997.39999999999
$ws.Range("N140").Value = -99357.39999999999

# WVR!row62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 636242.3
$ws.Range("I62").Value = 761600.8
$ws.Range("J62").Value = 9450
$ws.Range("K62").Value = 761600.8
$ws.Range("L62").Value = 9450
$ws.Range("M62").Value = -760976.8
$ws.Range("N62").Value = -10698

# WVR!row65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 636242.3
$ws.Range("I65").Value = 761600.8
$ws.Range("J65").Value = 9450
$ws.Range("K65").Value = 3808004
$ws.Range("L65").Value = 47250
$ws.Range("M65").Value = -3804884
$ws.Range("N65").Value = -53490

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 968535
$ws.Range("I136").Value = 1931324.8
$ws.Range("K136").Value = 5793974.4
$ws.Range("M136").Value = -5791424.4
